# "Issues status is updated"
# - Bug sheet: mark issue #1 (row 2) as Closed instead of Open.
# - Bug sheet: log a new issue (#9, row 11) describing the invoice
#   line-item Dto.Convert() null-reference bug, status "Check in",
#   remark "Opne".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug")

# Row 2 (issue #1): Status Open -> Closed
$ws.Range("D2").Value = "Closed"

# New row 11 (issue #9)
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Vanilla.Accontat.Facade.Invoice.LineItem.Server.cs -> public override BinAff.Facade.Library.Dto Convert(BinAff.Core.Data data) has thrown object reference not set to an reference error while loading old record"
$ws.Range("C11").Value = "Check in"
$ws.Range("D11").Value = "Opne"

$ws.Rows.Item(11).RowHeight = 45

$ws.Range("D2").Select()
